$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24. This shifts the existing rows 24-137
# down to 25-138 (carrying all of their column values along), exactly
# matching the "shift" pattern seen across the rest of the sheet.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new data point.
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 45243
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = "Tropicales y subtropicales"
$ws.Range("I24").Value = 100108004
$ws.Range("J24").Value = "Papaya"
$ws.Range("K24").Value = "Cultivar IV Región"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 90
$ws.Range("N24").Value = 26000
$ws.Range("O24").Value = 26000
$ws.Range("P24").Value = 26000
$ws.Range("Q24").Value = "`$/bandeja 10 kilos"
$ws.Range("R24").Value = "Provincia del Elquí"
$ws.Range("S24").Value = 2600
$ws.Range("T24").Value = 10
